$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new review row was inserted at row 23, pushing the existing rows 23-31
# down to 24-32. Use Excel's native row insert so formatting/styles carry
# over exactly like the rest of the data (style id 2 on column C).
$ws.Rows("23:23").Insert()

# Populate the newly inserted row with the new review's data (no comment
# text for this review, so column B is left blank).
$ws.Cells.Item(23, 1).Value = 5
$ws.Cells.Item(23, 3).Value = 46006.64905208333
$ws.Cells.Item(23, 4).Value = "YmIzNjc4MmUtMjk5Mi00NTY3LWE5ZTEtY2VkYTU4MWM3N2NlOjU3MDE2"
